# QUARTERLY_DATA_DOWNLOAD_ALL.xlsx refresh:
# - Column P ("Unemployment Rate") values are refreshed from the latest
#   source extract, which rounds the quarterly averages to 1 decimal place
#   instead of the old repeating-decimal figures.
# - N122 ("Employment level" for the most recent quarter) is newly
#   populated now that the source data for that quarter is available.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$unemploymentRateUpdates = @{
    14 = 8.8
    16 = 7.3
    17 = 6.7
    20 = 5.4
    24 = 4.2
    26 = 4
    28 = 4.4
    30 = 4.5
    31 = 4.6
    33 = 5.3
    34 = 4.8
    37 = 4.8
    38 = 5.1
    39 = 4.7
    42 = 4.4
    43 = 4.9
    45 = 4.5
    47 = 4.8
    49 = 4.6
    53 = 5.2
    54 = 5.4
    55 = 5.9
    56 = 7.4
    57 = 8.5
    58 = 10.8
    61 = 13.8
    63 = 14.4
    65 = 15.6
    66 = 15.2
    67 = 15.1
    69 = 15.9
    70 = 15.9
    71 = 15.5
    74 = 14.5
    75 = 14.4
    76 = 13.3
    77 = 12.9
    78 = 12.8
    80 = 11.7
    81 = 10.9
    82 = 10.6
    84 = 9.6
    85 = 9.4
    86 = 9
    88 = 8.2
    89 = 7.5
    90 = 7.3
    91 = 6.7
    92 = 6.7
    93 = 6.4
    94 = 6
    95 = 5.9
    99 = 5.2
    102 = 4.9
    103 = 5.1
    105 = 6.3
    106 = 7.3
    107 = 7
    108 = 5.5
    109 = 5.2
    110 = 4.9
    112 = 4.3
    113 = 4.4
    114 = 4.2
    115 = 4.2
    116 = 4.4
    117 = 4.5
    118 = 4.3
    119 = 4.3
    120 = 4.3
    121 = 4.3
    122 = 4
}

foreach ($row in $unemploymentRateUpdates.Keys) {
    $ws.Range("P$row").Value = $unemploymentRateUpdates[$row]
}

# New Employment level figure for row 122 (previously blank)
$ws.Range("N122").Value = 2794.1

